$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 4 (the first product row).
# This shifts rows 4..10 down to 5..11 and keeps merged cells / row heights
# of the shifted rows intact automatically.
$ws.Rows.Item(4).Insert()

# Copy formatting (fonts, fills, borders, number format, etc.) from the row
# that is now row 5 (originally row 4) onto the new, still-blank row 4, so
# that the new row matches the look of the rest of the table.
$ws.Range("A5:N5").Copy()
$ws.Range("A4:N4").PasteSpecial(-4122)

# Re-create the merged cell regions for the new row (merges are not carried
# over by PasteSpecial).
$ws.Range("B4:G4").Merge()
$ws.Range("H4:K4").Merge()
$ws.Range("L4:M4").Merge()

# Restore the exact per-row heights used by the source report for the
# product rows (they don't follow a simple shift pattern).
$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 24.75

# Populate the new first product row with the newly added drug.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "AVEROTHIAZIDE 5/20/12.5MG 30 F.C. TABS"
$ws.Range("H4").Value = "1:0"
$ws.Range("L4").Value = 93
$ws.Range("N4").Value = "1:0"

# Renumber the "م" (index) column for the rows that were pushed down, since
# they keep their original literal numbers after the insert.
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 6

# Update the running total (K column on the totals row, now row 10) to
# reflect the sum of the "الرصيد الحالي" (L) column across all product rows.
$ws.Range("K10").Value = 301

# The totals row grew slightly taller after the edit.
$ws.Rows.Item(10).RowHeight = 26.25
